# Fourth commit: completed rectangle.py as well as unit tests, and completed activity_02_main.py
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: __init__ - Attribute set to input values.
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "color = ""Red""`nlength = 8`nwidth = 4"
$ws.Range("G7").Value = "Object created with expected values"

# Row 8: __init__ - Exception raised when color is blank
$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "color = """"`nlength = 8`nwidth = 4"
$ws.Range("G8").Value = "ValueError"

# Row 9: __init__ - Exception raised when length is not an integer.
$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "color = ""Red""`nlength = ""hi""`nwidth = 4"
$ws.Range("G9").Value = "ValueError"

# Row 10: __init__ - Exception raised when width is not an integer.
$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "color = ""Red""`nlength = 8`nwidth = ""hi"""
$ws.Range("G10").Value = "ValueError"

# Row 11: __str__ - Returns string formatted appropriately
$ws.Range("E11").Value = "Rectangle = Rectangle(""Red"", 8, 4)"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "String returned in correct format"

# Row 12: calculate_area - Returns correct calculated value.
$ws.Range("E12").Value = "Rectangle = Rectangle(""Red"", 8, 4)"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "self.rectangle.area = 32"

# Row 13: calculate_perimeter - Returns correct calculated value.
$ws.Range("E13").Value = "Rectangle = Rectangle(""Red"", 8, 8)"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "self.rectangle.perimeter = 24"

# Selection / view state updates observed in the diff
$ws.Activate()
$ws.Range("G13").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
